$d = $word.ActiveDocument

# The page header carries the MarineGEO logo as an inline picture. This
# pass tightens up the header image sizing (part of a broader "updating
# headers and spacing" cleanup across the protocol docs): the picture is
# scaled down from its old extent to a slightly smaller one, keeping its
# aspect ratio.
#
#   old: cx=1156560 EMU (91.06771653543306 pt) x cy=395653 EMU (31.153779527559056 pt)
#   new: cx=1149366 EMU (90.50125984251969  pt) x cy=393192 EMU (30.96 pt)
#
# Word COM's InlineShape.Width / .Height are expressed in points;
# 1 point = 12700 EMU.
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$shp = $hdr.Range.InlineShapes.Item(1)

$shp.Width  = 1149366 / 12700
$shp.Height = 393192 / 12700
